# thesis/main_new_formulation/sets.xlsx
# Corretto errore in new formulation: identity_rcot deve essere definita come
# flows x activities2 (e non viceversa).
#
# _set_a2!A2:A16 (the "a2_names") are replaced with the real activity names
# taken from _set_a!A2:A16 (instead of the old "<name>2" placeholders).
# _set_a2!B2:B16 ("flows") are left untouched.
#
# Also tidies the leftover stray formatting column (was column E) by
# removing it and re-creating an (empty, formatted-only) column D, matching
# the "power_tril" cleanup mentioned in the commit message.

$wb = $excel.ActiveWorkbook

$wsA  = $wb.Worksheets.Item("_set_a")
$wsA2 = $wb.Worksheets.Item("_set_a2")

# --- _set_a2: column A becomes the activity names (flows x activities2) ---
for ($r = 2; $r -le 16; $r++) {
    $wsA2.Range("A$r").Value = $wsA.Range("A$r").Value2
}

# --- _set_a2: drop the old stray column E ---
$wsA2.Columns.Item(5).Delete()

# --- _set_a2: re-create the (empty) formatted column D, rows 2-15 ---
$fmtSrc = $wsA.Range("D9")
$fmtSrc.Copy()
for ($r = 2; $r -le 15; $r++) {
    $wsA2.Range("D$r").PasteSpecial(-4122)
}

# --- _set_a2: column B width + zoom + selection ---
$wsA2.Columns.Item(2).ColumnWidth = 15.83

$wsA2.Activate()
$excel.ActiveWindow.Zoom = 101
$wsA2.Range("D7").Select()

# --- _set_a: just the remembered selection moves to D1 ---
$wsA.Activate()
$wsA.Range("D1").Select()

# restore _set_a2 as the active sheet (it was tabSelected in the original file)
$wsA2.Activate()
